$d = $word.ActiveDocument

# --- Paragraph 1: update date and title (two runs split by a line break) ---
$d.Content.Find.Execute("המאמר היומי של מייק - 09.04.25", $false, $false, $false, $false, $false, $true, 1, $false, "המאמר היומי של מייק - 07.04.25", 2) | Out-Null
$d.Content.Find.Execute("O1-CODER: AN O1 REPLICATION FOR CODING", $false, $false, $false, $false, $false, $true, 1, $false, "JETFORMER: AN AUTOREGRESSIVE GENERATIVE MODEL OF RAW IMAGES AND TEXT", 2) | Out-Null

# --- Paragraphs 2-6: replace body text (use Range delete+insert to avoid Find/Replace AutoCorrect mangling quotes) ---
$r2 = $d.Paragraphs.Item(2).Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Delete()
$r2.InsertAfter("הסקירה של היום היא מאמר המשך(למרות שאין כאן VAE) של סקירתי האחרונה (מ 05.04.25). המאמר שנסקור היום מציע שיטה לאימון מודל מולטימודלי כאשר מודל אוטורגרסיבי אחד מאומן לשתי המודליות (תמונות וטקסט יחד). ברוב המודלים המולטי-מודליים יש אנקודרים שונים לטקסט ותמונה ולדעת מחברי המאמר זה עלול להוות בעיה (אני סוג של מבין את זה). אז המאמר מציע לאמן טרנספורמר אוטורגרסיבי לשתי המודליות יחד.")
$r3 = $d.Paragraphs.Item(3).Range
$r3.MoveEnd(1, -1) | Out-Null
$r3.Delete()
$r3.InsertAfter("אז איך הדבר הזה עובד בעצם? המאמר מציע להשתמש במודל מאומן של זרימה מנורמלת (Normalized Flows or NF) לבניית ייצוג התמונה. מודל NF מאמן מיפוי הפיך ולכן lossless ממרחב הדאטה (תמונה) למרחב בעל התפלגות פשוטה (נגיד גאוסית סטנדרטית). בדרך כלל מיפוי זה נבנה על ידי הרכבה (composition) של כמה מיפוים פשוטים (נגיד על תת-קבוצה קטנה של מימדים) וכל המיפויים האלו מאמנים יחד כאשר המטרה היא למקסם את הנראות (likelihood) של הדאטה תחת המיפוי הזה. למעשה המחברים מאמנים NF עבור כל פאץ' בתמונה (ייצוג פאץ' נקרא טוקן ויזואלי). ")
$r4 = $d.Paragraphs.Item(4).Range
$r4.MoveEnd(1, -1) | Out-Null
$r4.Delete()
$r4.InsertAfter("אז המחברים מאמנים יחד מודל NF לייצוג תמונה יחד עם טרנספורמר אוטורגרסיבי לגנרוט תמונה וטקסט. כלומר בהינתן תיאור התמונה והתמונה עצמה (הסדר בהזנה של פיסות דאטה חשוב!) הטרנספורמר אוטורגרסיבי מאומן לפלוט את ייצוגי הטוקנים הויזואליים אחרי NF (שמאומנים יחד עם הטרנספורמר). כאשר תמונה מוזנת לפני התיאור שלה הטרנספורמר מאומן לשחזר את ייצוג הטוקנים הטקסטואליים. כמו בסקירה הקודמת (GIVT) המודל חוזה פרמטרים של ה-gaussian mixture עבור כל טוקן והייצוג נדגם משם.")
$r5 = $d.Paragraphs.Item(5).Range
$r5.MoveEnd(1, -1) | Out-Null
$r5.Delete()
$r5.InsertAfter("המאמר גם מציע להעלות את הרובסטיות של ייצוגים המופקים על ידי המודל האוטורגרסיבי המאומן עם הרעשת דאטה(רק דאטה ויזואלי מורעש לפי הבנתי) מדורגת (סוג של למידת curriculum). בהתחלה מוסיפים לדאטה רעש חזק יותר כך שהמודל אוכל ללמוד את הפרטים ״הגסים״ של הדאטה ומורידים אותו במהלך האימון כך שהמודל ילמד גם את הפרטים העדינים יותר של הדאטה.")
$r6 = $d.Paragraphs.Item(6).Range
$r6.MoveEnd(1, -1) | Out-Null
$r6.Delete()
$r6.InsertAfter("https://arxiv.org/abs/2411.19722")

# --- Remove trailing paragraphs 9, 8, 7 (delete from the end to keep indices stable) ---
$d.Paragraphs.Item(9).Range.Delete()
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
